# Append the August 24th, 2020 row to the "out_vars" sheet.
# Source data (raw + clean) mirrors the row pattern already present for
# prior dates: Fecha, Confirmados, Negativos, Sospechosos, Defunciones,
# Porcentaje hospitalizados.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

$row = 86

# Force column A to be written as plain text (matching every other date
# cell in the column, e.g. "2020-08-23") instead of letting Excel's
# autodetection turn it into a date serial number, then drop the
# number-format override so the cell keeps the sheet's default style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-08-24"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 563705
$ws.Cells.Item($row, 3).Value = 622932
$ws.Cells.Item($row, 4).Value = 77198
$ws.Cells.Item($row, 5).Value = 60800
$ws.Cells.Item($row, 6).Value = 25.81
